# This script applies the commit's edits:
#  1) Collapses several runs that had been split apart by Word's spell-
#     checker (w:proofErr spellStart/spellEnd wrapping a "misspelled" word)
#     back into a single contiguous run - a purely cosmetic/internal change
#     that leaves the visible text identical. We do this by Find&Replace'ing
#     the *whole* phrase (which spans the split runs) with itself; Word's
#     find/replace re-literalizes the matched range into one run and drops
#     the proofErr bookmarks that fell inside it.
#  2) Fixes the "coeffient" -> "coefficient" typo.
#  3) Inserts a new italicized aside into the final answer paragraph.

$d = $word.ActiveDocument

function Replace-AllWhole([string]$findText, [string]$replaceText) {
    # Replace every occurrence of $findText with $replaceText, walking
    # forward through the document so all matches (not just the first) are
    # handled.
    $rng = $d.Content
    $rng.Start = 0
    while ($rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)) {
        $rng.Collapse(0)
        $rng.End = $d.Content.End
    }
}

# --- 1) Text-neutral run merges (drop w:proofErr, rejoin runs) ---------

Replace-AllWhole "runfile('/Users/zli/Desktop/Multiple-Linear-Regression/multiple_linear_regression.py', wdir='/Users/zli/Desktop/Multiple-Linear-Regression')" "runfile('/Users/zli/Desktop/Multiple-Linear-Regression/multiple_linear_regression.py', wdir='/Users/zli/Desktop/Multiple-Linear-Regression')"

Replace-AllWhole "Covariance Type:            nonrobust                                         " "Covariance Type:            nonrobust                                         "

Replace-AllWhole "                      coef    std err          t      P>|t|      [0.025      0.975]" "                      coef    std err          t      P>|t|      [0.025      0.975]"

Replace-AllWhole "                     coef    std err          t      P>|t|      [0.025      0.975]" "                     coef    std err          t      P>|t|      [0.025      0.975]"

Replace-AllWhole "Prob(Omnibus):                  0.000   Jarque-Bera (JB):               23.231" "Prob(Omnibus):                  0.000   Jarque-Bera (JB):               23.231"

Replace-AllWhole "Prob(Omnibus):                  0.001   Jarque-Bera (JB):               21.150" "Prob(Omnibus):                  0.001   Jarque-Bera (JB):               21.150"

Replace-AllWhole "Prob(Omnibus):                  0.940   Jarque-Bera (JB):                0.070" "Prob(Omnibus):                  0.940   Jarque-Bera (JB):                0.070"

Replace-AllWhole "Prob(Omnibus):                  0.947   Jarque-Bera (JB):                0.140" "Prob(Omnibus):                  0.947   Jarque-Bera (JB):                0.140"

Replace-AllWhole "The estimated parameter won't change but the confidence interval (or the c.i. range) could shrink by approx sqrt(2)." "The estimated parameter won't change but the confidence interval (or the c.i. range) could shrink by approx sqrt(2)."

Replace-AllWhole "2 features -> desicion tree/boosting/deep learning is not adequate." "2 features -> desicion tree/boosting/deep learning is not adequate."

# "Anormaly Dection (to be reviewed):" starts its paragraph with the
# w:proofErr straight away (no preceding run in that paragraph), so a plain
# whole-phrase replace can't scoop it up - extend the search to start one
# character earlier (the paragraph mark ending the previous, empty,
# paragraph) so the match range brackets the proofErr and Word's
# find/replace discards it along with the run split.
Replace-AllWhole "`rAnormaly Dection (to be reviewed):" "`rAnormaly Dection (to be reviewed):"

# --- 2) Typo fix: coeffient -> coefficient ------------------------------

# The w:proofErr spellEnd sits immediately after "coeffient" and before the
# following run's text, so extend the match one character past the word to
# pull that closing proofErr into the replaced range too.
Replace-AllWhole "which means the coeffient " "which means the coefficient "

# --- 3) Insert italicized aside into the final answer paragraph --------

# Locate the end of "...or 10.5%" (right before the existing period) and
# splice in the new aside there, so the original period ends up right after
# the new "10%)" text: "...10.5% (alternatively ... or 10%). The Z score..."
$insertPoint = $d.Content
$insertPoint.Start = 0
$insertPoint.Find.Execute("fraudulence of exp(0.10)~=1.105 or 10.5%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint.Collapse(0)
$insertPoint.InsertAfter(" (alternatively the increase in the log-odds of fraud of 0.1 or 10%)")

# Now italicize the whole clause from "an increase in the odds" through the
# inserted "... or 10%)." (matching the source, the trailing period of that
# clause is italicized too).
$italicRng = $d.Content
$italicRng.Start = 0
$italicRng.Find.ClearFormatting()
$ok = $italicRng.Find.Execute("an increase in the odds of fraudulence of exp(0.10)~=1.105 or 10.5% (alternatively the increase in the log-odds of fraud of 0.1 or 10%).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    # NOTE: Font.ItalicBi is intentionally *not* used here - in this
    # runtime it mis-scopes and stamps <w:iCs/> on every run in the whole
    # paragraph (not just the matched range). Font.Italic alone correctly
    # scopes to the matched range and is sufficient for the visible italics.
    $italicRng.Font.Italic = $true
}
